# "Generate Report for Handback"
#
# The 658c05b7-... record (row 3) failed its handback transform instead of
# being ready for handoff, so its status changes everywhere it is shown
# (Overview!B3/C3, zh-cn!C3, de-de!C3) and each language sheet gets a new
# Error Detail (column K) explaining the mismatched file name.

$wb = $excel.ActiveWorkbook

# Status text used on the Overview sheet and on each language sheet's
# "Status" column (C) for the 658c05b7-... row is now "Handback transform
# failed" instead of "Ready for handoff". Replace() updates every
# occurrence across all worksheets in one pass.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "Handback transform failed")
}

# Record the handback/handoff file-name mismatch in the "Error Detail"
# column (K) of row 3 for each language sheet.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K3").Value = "Handback file name: bvrf2kvn.spx is different with handoff file name: 658c05b7-2ba7-4441-a45d-bf993581e007.ec8f02086e8621d446274e4ad4674fd819ae9e22.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "Handback file name: bvrf2kvn.spx is different with handoff file name: 658c05b7-2ba7-4441-a45d-bf993581e007.ec8f02086e8621d446274e4ad4674fd819ae9e22.de-de."
